# all work for 2/24/2023 - Issues #70 #66 #27
#
# The "Computer Programming" row (originally row 2) moves to the bottom of
# the table (new row 6). The stray surrounding quotes on the video-url
# value ("/videos/saiareact.mp4") are stripped everywhere that value is
# still used, and the Networking row's video-url (C5) is cleared out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Copy the original "Computer Programming" row (row 2, with its
#    formatting/styles) down to the new last row (row 6).
$ws.Range("A2:D2").Copy($ws.Range("A6:D6"))

# Fix up the video endpoint in the freshly copied row: strip the stray
# literal quote characters that were baked into the old string.
$ws.Range("C6").Value = "/videos/saiareact.mp4"

# Match the tall "wrap text" row height used by the rest of the table.
$ws.Rows(6).RowHeight = 409.5

# 2) Blank out row 2 (keep the cell formatting/styles, drop the values).
$ws.Range("A2").Value = ""
$ws.Range("B2:D2").ClearContents()
$ws.Rows(2).AutoFit()

# 3) Strip the stray quotes from the video endpoint still referenced by
#    the Electronics/IT Support (row 3) and Information Systems
#    Management (row 4) rows.
$ws.Range("C3").Value = "/videos/saiareact.mp4"
$ws.Range("C4").Value = "/videos/saiareact.mp4"

# 4) The Networking and Cyber Security row no longer has a video endpoint.
$ws.Range("C5").ClearContents()

# 5) Leave the selection where the author left it.
$ws.Range("C5").Select()
